# Commit: "changed font to Verdana"
#
# The legend rectangles (shapes 1-9) and the accompanying key table on
# slide 1 had their run text shrunk from 16pt to 12pt and switched from
# the theme's default font to Verdana (both the Latin and the
# complex-script typeface slots).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Set-VerdanaFont {
    param($textRange)
    $textRange.Font.Size = 12
    $textRange.Font.Name = "Verdana"
    $textRange.Font.NameComplexScript = "Verdana"
}

# --- The 9 colour-coded legend rectangles ---
for ($i = 1; $i -le 9; $i++) {
    $shp = $s.Shapes.Item($i)
    Set-VerdanaFont $shp.TextFrame.TextRange
}

# --- The "Code"/"Phylum" key table ---
$tblShape = $s.Shapes.Item(10)
$tbl = $tblShape.Table

for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
        $cell = $tbl.Cell($r, $c)
        Set-VerdanaFont $cell.Shape.TextFrame.TextRange
    }
}

# Shrinking the text lets PowerPoint re-fit the table to its rows'
# specified heights, which is what actually produced the shorter
# graphicFrame extent in the saved file.
for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    $tbl.Rows($r).Height = 24.066692913385825
}
